# Add "Expected Result" and "Actual Result" columns to the API test sheet,
# between PARAM (C) and TEST RESULT (D), shifting the old D/E columns to F/G.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("API")

# Insert two new blank columns at D; old D (TEST RESULT) and E (Comment)
# shift right to F and G.
$ws.Columns.Item(4).Insert()
$ws.Columns.Item(4).Insert()

# New header cells.
$ws.Range("D1").Value = "Expected Result"
$ws.Range("E1").Value = "Actual Result"

# New row 2 data (showreg.php test case).
$ws.Range("D2").Value = "It should show all records of registration form in JSON (Java Script Object Notation) format"
$ws.Range("E2").Value = "as per expected"

# Row 3 (reg.php test case) - expanded PARAM list plus new expected/actual/comment text.
$ws.Range("C3").Value = "email, pwd,name,dob,mobile"
$ws.Range("D3").Value = "It should add new record into database  and show status 1 and done user added message"
$ws.Range("E3").Value = "as per expected"
$ws.Range("G3").Value = "blank data inserted"

# Column widths for the two new columns, and restore the previous widths on
# the shifted-right columns (F keeps the old D width, G keeps the old E width).
$ws.Range("D1").ColumnWidth = 81.66666666666667
$ws.Range("E1").ColumnWidth = 81.66666666666667
$ws.Range("F1").ColumnWidth = 22.166666666666668
$ws.Range("G1").ColumnWidth = 69.66666666666667

# Move the active selection to A9, as in the saved workbook.
$ws.Range("A9").Select()
